# Commit: "Successfull - Commented GitLatch Commit @ 2022-12-28-8-8-35-710"
#
# Two changes were made in the canonical OOXML:
#
# 1) ppt/presentation.xml gained an empty p15:sldGuideLst extension block
#    (no actual guides were added - the list is empty). This corresponds to
#    the presentation-level Guides feature being touched in the UI.
# 2) The "Picture Placeholder 2" shape on the "Picture with Caption" slide
#    layout (ppt/slideLayouts/slideLayout9.xml) had its empty paragraph
#    (<a:endParaRPr/> only) replaced with a run containing the placeholder
#    prompt text "Click icon to add picture".

$p = $ppt.ActivePresentation

# --- 1) Touch the presentation-level Guides so PowerPoint (re)writes the
#        (possibly empty) slide-guide list extension on save. ---
$guides = $p.Guides
$newGuide = $guides.Add(1, 100)
if ($newGuide) {
    $newGuide.Delete()
}
$ppt.DisplayGuides = $true

# --- 2) Set the picture placeholder's prompt text on the "Picture with
#        Caption" layout (the 9th custom layout on the slide master). ---
$master = $p.SlideMaster
$customLayouts = $master.CustomLayouts
for ($i = 1; $i -le $customLayouts.Count; $i++) {
    $layout = $customLayouts.Item($i)
    if ($layout.Name -eq "Picture with Caption") {
        for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
            $shape = $layout.Shapes.Item($j)
            if ($shape.Name -eq "Picture Placeholder 2") {
                $shape.TextFrame.TextRange.Text = "Click icon to add picture"
            }
        }
    }
}
